# This workbook is a weekly price-report extract. The update inserts one
# new week's worth of data (2 rows: "Primera" and "Segunda" quality grades)
# at the top of the "Acelga" price history block (which starts at row 357),
# pushing all the older rows down by two positions. The two oldest rows
# that fall off the bottom of the original range become the new last rows
# of the sheet (485 and 486).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 357, shifting
# rows 357:484 down to 359:486 (dimension grows from R484 to R486
# automatically).
$ws.Rows.Item(357).Resize(2).Insert()

# Row 357 - "Primera" grade, new week
$ws.Cells.Item(357, 1).Value = 8
$ws.Cells.Item(357, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(357, 3).Value = "Coquimbo"
$ws.Cells.Item(357, 4).Value = 44809
$ws.Cells.Item(357, 5).Value = 4
$ws.Cells.Item(357, 6).Value = 100112009
$ws.Cells.Item(357, 7).Value = "Acelga"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 2540
$ws.Cells.Item(357, 11).Value = 650
$ws.Cells.Item(357, 12).Value = 700
$ws.Cells.Item(357, 13).Value = 675
$ws.Cells.Item(357, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(357, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(357, 16).Value = 338
$ws.Cells.Item(357, 17).Value = 2
$ws.Cells.Item(357, 18).Value = "Hortaliza"

# Row 358 - "Segunda" grade, new week
$ws.Cells.Item(358, 1).Value = 8
$ws.Cells.Item(358, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(358, 3).Value = "Coquimbo"
$ws.Cells.Item(358, 4).Value = 44809
$ws.Cells.Item(358, 5).Value = 4
$ws.Cells.Item(358, 6).Value = 100112009
$ws.Cells.Item(358, 7).Value = "Acelga"
$ws.Cells.Item(358, 8).Value = "Sin especificar"
$ws.Cells.Item(358, 9).Value = "Segunda"
$ws.Cells.Item(358, 10).Value = 1400
$ws.Cells.Item(358, 11).Value = 550
$ws.Cells.Item(358, 12).Value = 600
$ws.Cells.Item(358, 13).Value = 575
$ws.Cells.Item(358, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(358, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(358, 16).Value = 288
$ws.Cells.Item(358, 17).Value = 2
$ws.Cells.Item(358, 18).Value = "Hortaliza"
